# Adds two new floraväkteri observation rows (12 and 13) to the 'Artfynd' sheet,
# replicating a new data export batch appended to the existing table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- helpers -------------------------------------------------------------

# Plain text that Excel would not misinterpret as a number/date/bool.
function Set-Text($row, $col, $value) {
    $ws.Cells.Item($row, $col).Value = $value
}

# Text that *looks* numeric/date-like (e.g. "1", "2023-10-02") - a leading
# apostrophe forces Excel to store it as text instead of auto-converting it;
# the apostrophe itself is not part of the stored value. Style is reset back
# to Normal afterwards so the quote-prefix flag doesn't linger on the cell.
function Set-ForcedText($row, $col, $value) {
    $ws.Cells.Item($row, $col).Value = "'" + $value
    $ws.Cells.Item($row, $col).Style = "Normal"
}

# An explicitly-present-but-blank text cell (source export emits these for
# optional fields that are present in the schema but unset for this record).
function Set-EmptyText($row, $col) {
    $ws.Cells.Item($row, $col).Value = "'"
    $ws.Cells.Item($row, $col).Style = "Normal"
}

function Set-Number($row, $col, $value) {
    $ws.Cells.Item($row, $col).Value = $value
}

function Set-Bool($row, $col, $value) {
    $ws.Cells.Item($row, $col).Value = $value
}

# --- row data --------------------------------------------------------------

# Row 12
Set-Number 12 1 112488501
Set-Number 12 2 90155
Set-Text 12 3 'Ovaliderad'
Set-Text 12 4 'LC'
Set-Number 12 5 6031
Set-Text 12 6 'Blomkålssvamp'
Set-Text 12 7 'Sparassis crispa'
Set-Text 12 8 '(Wulfen:Fr.) Fr.'
Set-ForcedText 12 9 '1'
Set-Text 12 10 'fruktkroppar'
Set-EmptyText 12 11
Set-EmptyText 12 14
Set-Text 12 16 'Hyltåkra, Sm'
Set-Number 12 17 428943
Set-Number 12 18 6274121
Set-Number 12 19 5
Set-Text 12 20 'Kronoberg'
Set-Text 12 21 'Ljungby'
Set-Text 12 22 'Småland'
Set-Text 12 23 'Hamneda'
Set-ForcedText 12 25 '2023-10-02'
Set-ForcedText 12 27 '2023-10-02'
Set-Text 12 29 'Rullstensås. Vid tall.'
Set-Bool 12 30 $false
Set-Bool 12 31 $false
Set-EmptyText 12 32
Set-Bool 12 33 $false
Set-EmptyText 12 46
Set-Text 12 49 'Krister Wahlström'
Set-Text 12 50 'Krister Wahlström'
Set-EmptyText 12 51

# Row 13
Set-Number 13 1 112486878
Set-Number 13 2 96720
Set-Text 13 3 'Ovaliderad'
Set-Text 13 4 'VU'
Set-Number 13 5 220787
Set-Text 13 6 'Knärot'
Set-Text 13 7 'Goodyera repens'
Set-Text 13 8 '(L.) R. Br.'
Set-ForcedText 13 9 '101'
Set-Text 13 10 'stjälkar/strån/skott'
Set-Text 13 11 'blomning'
Set-EmptyText 13 12
Set-EmptyText 13 14
Set-Text 13 16 'Hyltåkra, Sm'
Set-Number 13 17 428966
Set-Number 13 18 6274117
Set-Number 13 19 61
Set-Text 13 20 'Kronoberg'
Set-Text 13 21 'Ljungby'
Set-Text 13 22 'Småland'
Set-Text 13 23 'Hamneda'
Set-Text 13 24 'G-Lju-1388'
Set-ForcedText 13 25 '2023-10-02'
Set-ForcedText 13 27 '2023-10-02'
Set-Text 13 29 '1 blommsamling. Skogen orörd.'
Set-Bool 13 30 $false
Set-Bool 13 31 $false
Set-EmptyText 13 32
Set-Bool 13 33 $false
Set-EmptyText 13 46
Set-Text 13 49 'Krister Wahlström'
Set-Text 13 50 'Krister Wahlström'
Set-Text 13 51 'Floraväkteri Sverige'

